$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.299.17"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "2.586.74"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'517.76"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("E6").Value = "  -2.37%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -1.43%  "

$ws.Range("D9").Value = "2.605.39"
$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("D10").Value = "'6.46"
$ws.Range("E10").Value = "  -2.67%  "

$ws.Range("E11").Value = "  -3.02%  "

$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").Value = "3.039.46"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").Value = "58.263.43"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").Value = "'20.27"
$ws.Range("E16").Value = "  -3.52%  "

$ws.Range("E17").Value = "  -2.96%  "

$ws.Range("D18").Value = "2.560.89"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").Value = "'336.69"
$ws.Range("E19").Value = "  -1.42%  "

$ws.Range("E20").Value = "  -3.64%  "

$ws.Range("D21").Value = "'10.17"
$ws.Range("E21").Value = "  -3.81%  "

$ws.Range("D22").Value = "'6.39"
$ws.Range("E22").Value = "  +0.76%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").Value = "'65.19"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").Value = "'0.400"
$ws.Range("E26").Value = "  -4.41%  "

$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").Value = "2.685.76"
$ws.Range("E28").Value = "  -2.49%  "

$ws.Range("E29").Value = "  -3.43%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").Value = "0.0₃0733"
$ws.Range("E31").Value = "  -8.00%  "

$ws.Range("D32").Value = "'6.09"
$ws.Range("E32").Value = "  -5.96%  "

$ws.Range("E33").Value = "  -2.90%  "

$ws.Range("D34").Value = "'18.66"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").Value = "'149.51"
$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").Value = "'3.96"
$ws.Range("E36").Value = "  -6.59%  "

$ws.Range("E37").Value = "  -6.33%  "

$ws.Range("D38").Value = "'0.849"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'36.20"
$ws.Range("E39").Value = "  -1.03%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.46"
$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("D41").Value = "'0.821"
$ws.Range("E41").Value = "  -11.87%  "

$ws.Range("E42").Value = "  -4.37%  "

$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").Value = "'272.84"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("D47").Value = "'0.0947"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").Value = "'18.66"
$ws.Range("E48").Value = "  -3.79%  "

$ws.Range("D49").Value = "'0.0518"
$ws.Range("E49").Value = "  -3.88%  "

$ws.Range("D50").Value = "1.966.79"
$ws.Range("E50").Value = "  -3.86%  "

$ws.Range("D51").Value = "'4.56"
$ws.Range("E51").Value = "  -4.95%  "
